$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 19-46: refreshed charger/terminal records (station name, terminal name, last-charge-end datetime)
$ws.Range("A19").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B19").Value = "103号直流"
$ws.Range("C19").Value = 45941.259837962964
$ws.Range("A20").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B20").Value = "404号直流"
$ws.Range("C20").Value = 45941.277685185189
$ws.Range("A21").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B21").Value = "603号直流"
$ws.Range("C21").Value = 45942.459050925929
$ws.Range("A22").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B22").Value = "406号直流"
$ws.Range("C22").Value = 45943.020914351851
$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "604号直流"
$ws.Range("C23").Value = 45943.03466435185
$ws.Range("A24").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B24").Value = "B01号直流"
$ws.Range("C24").Value = 45943.623032407406
$ws.Range("A25").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B25").Value = "801号直流"
$ws.Range("C25").Value = 45943.990763888891
$ws.Range("A26").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B26").Value = "201号直流"
$ws.Range("C26").Value = 45944.074282407404
$ws.Range("A27").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B27").Value = "203号直流"
$ws.Range("C27").Value = 45944.228055555555
$ws.Range("A28").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B28").Value = "704号直流"
$ws.Range("C28").Value = 45944.253993055558
$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "904号直流"
$ws.Range("C29").Value = 45944.263796296298
$ws.Range("A30").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B30").Value = "602号直流"
$ws.Range("C30").Value = 45944.51699074074
$ws.Range("A31").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B31").Value = "504号直流"
$ws.Range("C31").Value = 45944.527499999997
$ws.Range("A32").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B32").Value = "108号直流"
$ws.Range("C32").Value = 45944.530509259261
$ws.Range("A33").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B33").Value = "502号直流"
$ws.Range("C33").Value = 45944.54005787037
$ws.Range("A34").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B34").Value = "401号直流"
$ws.Range("C34").Value = 45944.554016203707
$ws.Range("A35").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B35").Value = "112号直流"
$ws.Range("C35").Value = 45944.583692129629
$ws.Range("A36").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B36").Value = "202号直流"
$ws.Range("C36").Value = 45944.60528935185
$ws.Range("A37").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B37").Value = "604号直流"
$ws.Range("C37").Value = 45944.616377314815
$ws.Range("A38").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B38").Value = "005B号直流"
$ws.Range("C38").Value = 45944.617847222224
$ws.Range("A39").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B39").Value = "008B号直流"
$ws.Range("C39").Value = 45944.62395833333
$ws.Range("A40").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B40").Value = "105号直流"
$ws.Range("C40").Value = 45944.624421296299
$ws.Range("A41").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B41").Value = "003B号直流"
$ws.Range("C41").Value = 45944.642442129632
$ws.Range("A42").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B42").Value = "006A号直流"
$ws.Range("C42").Value = 45944.650648148148
$ws.Range("A43").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B43").Value = "B03号直流"
$ws.Range("C43").Value = 45944.651828703703
$ws.Range("A44").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B44").Value = "306号直流"
$ws.Range("C44").Value = 45944.674907407411
$ws.Range("A45").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B45").Value = "203号直流"
$ws.Range("C45").Value = 45944.697951388887
$ws.Range("A46").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B46").Value = "304号直流"
$ws.Range("C46").Value = 45944.771261574075

# Rows 47-56: no longer have data, clear them out (keep styles/formatting intact)
$ws.Range("A47:C56").ClearContents()

# Update the saved selection/active cell
$ws.Range("D21").Select()
